$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 9261834
$ws.Range("I62").Value = 11113201
$ws.Range("K62").Value = 11113201
$ws.Range("M62").Value = -11112577

# Row 65
$ws.Range("H65").Value = 9261834
$ws.Range("I65").Value = 11113201
$ws.Range("K65").Value = 55566005
$ws.Range("M65").Value = -55562885

# Row 132
$ws.Range("H132").Value = 1209.1892
$ws.Range("I132").Value = 1084.6875
$ws.Range("J132").Value = 2006
$ws.Range("K132").Value = 3254.0625
$ws.Range("L132").Value = 6018
$ws.Range("M132").Value = -724.0625
$ws.Range("N132").Value = -11078

# Row 138
$ws.Range("H138").Value = 3257.2122
$ws.Range("I138").Value = 1613.0667
$ws.Range("J138").Value = 4627.3335
$ws.Range("K138").Value = 4839.2001
$ws.Range("L138").Value = 13882.0005
$ws.Range("M138").Value = 300.7999
$ws.Range("N138").Value = -24162.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 2098.9092
$ws.Range("I28").Value = 2098.9092
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2098.9092
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -1906.9092
$ws.Range("N28").ClearContents()

# Row 32
$ws.Range("H32").Value = 14619.263
$ws.Range("I32").Value = 11068.654
$ws.Range("J32").Value = 47166.5
$ws.Range("K32").Value = 11068.654
$ws.Range("L32").Value = 47166.5
$ws.Range("M32").Value = -10781.654
$ws.Range("N32").Value = -47740.5

# Row 56
$ws.Range("H56").Value = 339666.34
$ws.Range("J56").Value = 999999
$ws.Range("L56").Value = 999999
$ws.Range("N56").Value = -1001483

# Row 61
$ws.Range("H61").Value = 2950.2307
$ws.Range("I61").Value = 1863.4615
$ws.Range("J61").Value = 5123.769
$ws.Range("K61").Value = 1863.4615
$ws.Range("L61").Value = 5123.769
$ws.Range("M61").Value = -1651.4615
$ws.Range("N61").Value = -5547.769

# Row 99
$ws.Range("H99").Value = 2098.9092
$ws.Range("I99").Value = 2098.9092
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2098.9092
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 896.0907999999999
$ws.Range("N99").ClearContents()

# Row 110
$ws.Range("H110").Value = 1863.1
$ws.Range("I110").Value = 1863.1
$ws.Range("K110").Value = 1863.1
$ws.Range("M110").Value = 181.9000000000001

# Row 122
$ws.Range("H122").Value = 3040.7908
$ws.Range("I122").Value = 2291.6428
$ws.Range("K122").Value = 6874.928400000001
$ws.Range("M122").Value = -4424.928400000001

# Row 132
$ws.Range("H132").Value = 4438.3335
$ws.Range("I132").Value = 4210.0938
$ws.Range("K132").Value = 12630.2814
$ws.Range("M132").Value = -10100.2814

# Row 136
$ws.Range("H136").Value = 2950.2307
$ws.Range("I136").Value = 1863.4615
$ws.Range("J136").Value = 5123.769
$ws.Range("K136").Value = 5590.3845
$ws.Range("L136").Value = 15371.307
$ws.Range("M136").Value = -3040.3845
$ws.Range("N136").Value = -20471.307

# Row 140
$ws.Range("H140").Value = 76000
$ws.Range("J140").Value = 76000
$ws.Range("L140").Value = 76000
$ws.Range("N140").Value = -86360

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1940.619
$ws.Range("I20").Value = 1599.3636
$ws.Range("J20").Value = 2316
$ws.Range("K20").Value = 1599.3636
$ws.Range("L20").Value = 2316
$ws.Range("M20").Value = -1352.3636
$ws.Range("N20").Value = -2810

# Row 22
$ws.Range("H22").Value = 2931.125
$ws.Range("I22").Value = 4459.8
$ws.Range("K22").Value = 4459.8
$ws.Range("M22").Value = -4286.8

# Row 94
$ws.Range("H94").Value = 16670466
$ws.Range("I94").Value = 7695091.5
$ws.Range("K94").Value = 7695091.5
$ws.Range("M94").Value = -7694640.5

# Row 100
$ws.Range("H100").Value = 47000.5
$ws.Range("J100").Value = 47000.5
$ws.Range("L100").Value = 47000.5
$ws.Range("N100").Value = -49164.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 700
$ws.Range("I16").Value = 691.6667
$ws.Range("J16").Value = 750
$ws.Range("K16").Value = 691.6667
$ws.Range("L16").Value = 750
$ws.Range("M16").Value = -404.6667
$ws.Range("N16").Value = -1324

# Row 113
$ws.Range("H113").Value = 700
$ws.Range("I113").Value = 691.6667
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 691.6667
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = 1478.3333
$ws.Range("N113").Value = -5090

$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 1900
$ws.Range("J137").Value = 1900
$ws.Range("L137").Value = 5700
$ws.Range("N137").Value = -15900

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 3090550.2
$ws.Range("I7").Value = 8802
$ws.Range("J7").Value = 3370709
$ws.Range("K7").Value = 8802
$ws.Range("L7").Value = 3370709
$ws.Range("M7").Value = -8690
$ws.Range("N7").Value = -3370933

# Row 8
$ws.Range("H8").Value = 3090550.2
$ws.Range("I8").Value = 8802
$ws.Range("J8").Value = 3370709
$ws.Range("K8").Value = 8802
$ws.Range("L8").Value = 3370709
$ws.Range("M8").Value = -8663
$ws.Range("N8").Value = -3370987

# Row 113
$ws.Range("H113").Value = 4693.6
$ws.Range("I113").Value = 3129.625
$ws.Range("K113").Value = 3129.625
$ws.Range("M113").Value = -959.625

# Row 132
$ws.Range("H132").Value = 4410.886
$ws.Range("I132").Value = 3226.7917
$ws.Range("K132").Value = 9680.375100000001
$ws.Range("M132").Value = -7150.375100000001

# Row 135
$ws.Range("H135").Value = 59514.684
$ws.Range("J135").Value = 59514.684
$ws.Range("L135").Value = 59514.684
$ws.Range("N135").Value = -69654.68400000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1829.5
$ws.Range("I22").Value = 1505
$ws.Range("K22").Value = 1505
$ws.Range("M22").Value = -1210

# Row 27
$ws.Range("H27").Value = 1829.5
$ws.Range("I27").Value = 1505
$ws.Range("K27").Value = 1505
$ws.Range("M27").Value = -1398

# Row 122
$ws.Range("H122").Value = 5029.8125
$ws.Range("I122").Value = 3911.0435
$ws.Range("J122").Value = 7888.8887
$ws.Range("K122").Value = 11733.1305
$ws.Range("L122").Value = 23666.6661
$ws.Range("M122").Value = -9283.130500000001
$ws.Range("N122").Value = -28566.6661

# Row 136
$ws.Range("H136").Value = 5875.0835
$ws.Range("I136").Value = 3110.4
$ws.Range("J136").Value = 7849.857
$ws.Range("K136").Value = 9331.200000000001
$ws.Range("L136").Value = 23549.571
$ws.Range("M136").Value = -6781.200000000001
$ws.Range("N136").Value = -28649.571

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 6860.4375
$ws.Range("I136").Value = 6911.7
$ws.Range("J136").Value = 6775
$ws.Range("K136").Value = 20735.1
$ws.Range("L136").Value = 20325
$ws.Range("M136").Value = -18185.1
$ws.Range("N136").Value = -25425
